$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (P1) - shares the same style as the other date headers (O1, etc.)
$ws.Range("P1").Value = "29-jun"
$ws.Range("P1").NumberFormat = "@"

# New column data (P2:P11)
# P2 picks up a freshly-created style (general format + center alignment)
$ws.Range("P2").Value = 15
$ws.Range("P2").HorizontalAlignment = -4108

# P3:P11 reuse the existing "0" number format + center alignment style
$ws.Range("P3").Value = 16
$ws.Range("P3").HorizontalAlignment = -4108
$ws.Range("P3").NumberFormat = "0"

$ws.Range("P4").Value = 8
$ws.Range("P4").HorizontalAlignment = -4108
$ws.Range("P4").NumberFormat = "0"

$ws.Range("P5").Value = 12
$ws.Range("P5").HorizontalAlignment = -4108
$ws.Range("P5").NumberFormat = "0"

$ws.Range("P6").Value = 11
$ws.Range("P6").HorizontalAlignment = -4108
$ws.Range("P6").NumberFormat = "0"

$ws.Range("P7").Value = 20
$ws.Range("P7").HorizontalAlignment = -4108
$ws.Range("P7").NumberFormat = "0"

$ws.Range("P8").Value = 10
$ws.Range("P8").HorizontalAlignment = -4108
$ws.Range("P8").NumberFormat = "0"

$ws.Range("P9").Value = 18
$ws.Range("P9").HorizontalAlignment = -4108
$ws.Range("P9").NumberFormat = "0"

$ws.Range("P10").Value = 8
$ws.Range("P10").HorizontalAlignment = -4108
$ws.Range("P10").NumberFormat = "0"

$ws.Range("P11").Value = 13
$ws.Range("P11").HorizontalAlignment = -4108
$ws.Range("P11").NumberFormat = "0"

# Move the active selection to match the saved workbook state
$ws.Range("G4").Select() | Out-Null
